$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for c93c7021 row (row 3)
# This value is shared with de-de!H3 ("Correspond Handoff Datetime"), so both are updated.
$wsOverview.Range("G3").Value = "2016-08-20 14:48:54"
$wsDeDe.Range("H3").Value = "2016-08-20 14:48:54"

# zh-cn sheet: "Correspond Handoff Datetime" (H3) and "Correspond Handback DateTime" (K3) for c93c7021 row (row 3)
$wsZhCn.Range("H3").Value = "2016-08-20 14:48:50"
$wsZhCn.Range("K3").Value = "2016-08-20 14:49:12"

# de-de sheet: "Correspond Handback DateTime" (K3) for c93c7021 row (row 3)
$wsDeDe.Range("K3").Value = "2016-08-20 14:49:18"
